$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data between rows 21 and 22 (columns F:V). The "Indice"
#    (A) and "data_partida" (E) columns stay exactly where they were; only
#    the match/odds details move. Use a scratch row far outside the used
#    range as a temporary holding area, then clear it again.
# ---------------------------------------------------------------------------
$ws.Range("F21:V21").Copy($ws.Range("F1000:V1000"))
$ws.Range("F22:V22").Copy($ws.Range("F21:V21"))
$ws.Range("F1000:V1000").Copy($ws.Range("F22:V22"))
$ws.Range("F1000:V1000").ClearContents()

# ---------------------------------------------------------------------------
# 2) Append four new match rows (51-54) at the bottom of the table, copying
#    the formatting/styles from the last existing row (50) and then filling
#    in the new values.
# ---------------------------------------------------------------------------
$ws.Range("A50:V50").Copy($ws.Range("A51:V54"))

# Row 51: Brondby 2 - 3 FC Copenhagen
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "denmark"
$ws.Range("C51").Value = "superliga"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45193.58333333334
$ws.Range("F51").Value = "Brondby"
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = "FC Copenhagen"
$ws.Range("I51").Value = 3
$ws.Range("J51").Value = 2.74
$ws.Range("K51").Value = "17/09/2023 17:12"
$ws.Range("L51").Value = 2.89
$ws.Range("M51").Value = "24/09/2023 13:56"
$ws.Range("N51").Value = 3.54
$ws.Range("O51").Value = "17/09/2023 17:12"
$ws.Range("P51").Value = 3.52
$ws.Range("Q51").Value = "24/09/2023 13:56"
$ws.Range("R51").Value = 2.59
$ws.Range("S51").Value = "17/09/2023 17:12"
$ws.Range("T51").Value = 2.48
$ws.Range("U51").Value = "24/09/2023 13:56"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/denmark/superliga/brondby-fc-copenhagen/nmBPisfK/"

# Row 52: Silkeborg 2 - 0 Viborg
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "denmark"
$ws.Range("C52").Value = "superliga"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("E52").Value = 45193.58333333334
$ws.Range("F52").Value = "Silkeborg"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "Viborg"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1.98
$ws.Range("K52").Value = "18/09/2023 18:12"
$ws.Range("L52").Value = 1.95
$ws.Range("M52").Value = "24/09/2023 13:56"
$ws.Range("N52").Value = 3.73
$ws.Range("O52").Value = "18/09/2023 18:12"
$ws.Range("P52").Value = 3.84
$ws.Range("Q52").Value = "24/09/2023 13:56"
$ws.Range("R52").Value = 3.45
$ws.Range("S52").Value = "18/09/2023 18:12"
$ws.Range("T52").Value = 3.83
$ws.Range("U52").Value = "24/09/2023 13:56"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/denmark/superliga/silkeborg-viborg/0fi0ouXm/"

# Row 53: Randers FC 1 - 1 Aarhus
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "denmark"
$ws.Range("C53").Value = "superliga"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("E53").Value = 45193.66666666666
$ws.Range("F53").Value = "Randers FC"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "Aarhus"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 3.4
$ws.Range("K53").Value = "17/09/2023 17:12"
$ws.Range("L53").Value = 3.36
$ws.Range("M53").Value = "24/09/2023 15:57"
$ws.Range("N53").Value = 3.35
$ws.Range("O53").Value = "17/09/2023 17:12"
$ws.Range("P53").Value = 3.11
$ws.Range("Q53").Value = "24/09/2023 15:53"
$ws.Range("R53").Value = 2.22
$ws.Range("S53").Value = "17/09/2023 17:12"
$ws.Range("T53").Value = 2.41
$ws.Range("U53").Value = "24/09/2023 15:57"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/denmark/superliga/randers-fc-aarhus/voj4panf/"

# Row 54: Midtjylland 2 - 1 Odense
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "denmark"
$ws.Range("C54").Value = "superliga"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("E54").Value = 45193.75
$ws.Range("F54").Value = "Midtjylland"
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = "Odense"
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = 1.68
$ws.Range("K54").Value = "18/09/2023 18:12"
$ws.Range("L54").Value = 1.63
$ws.Range("M54").Value = "24/09/2023 17:59"
$ws.Range("N54").Value = 4.09
$ws.Range("O54").Value = "18/09/2023 18:12"
$ws.Range("P54").Value = 4.26
$ws.Range("Q54").Value = "24/09/2023 17:59"
$ws.Range("R54").Value = 4.85
$ws.Range("S54").Value = "18/09/2023 18:12"
$ws.Range("T54").Value = 5.35
$ws.Range("U54").Value = "24/09/2023 17:59"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/denmark/superliga/midtjylland-odense/UB0Uj19Q/"
